$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append a new row of mail-log data ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A4").Value = "Demo inplannen"
$logs.Range("B4").Value = "klantenservice@testbedrijf123.nl"
$logs.Range("C4").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Range("D4").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E4").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Range("F4").Value = "2025-08-13 20:01:43"
$logs.Range("G4").Value = "Nee"
$logs.Range("H4").Value = "Ja"
$logs.Range("I4").Value = "Nee"
$logs.Range("J4").Value = "Nee"

# Extend the conditional-formatting blocks so they keep covering the
# newly added row (row 4) for every formatted column.
$logs.Range("D2:D3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D4"))
$logs.Range("G2:G3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G4"))
$logs.Range("H2:H3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H4"))
$logs.Range("I2:I3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I4"))
$logs.Range("J2:J3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J4"))

# --- Sheet "Dashboard": bump the count for this category ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 3
